$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fuels")

# Row 17: syngas - ecoinvent
$ws.Range("A17").Value = "syngas - ecoinvent"
$ws.Range("C17").Formula = "=5.4/1.15"
$ws.Range("D17").Value = 0.927

# Row 18: syngas - PNNL
$ws.Range("A18").Value = "syngas - PNNL"
$ws.Range("B18").Value = 27.66299
$ws.Range("C18").Value = 27.66299
$ws.Range("D18").Value = 1.392744

# Row 19: wood chips - dry
$ws.Range("A19").Value = "wood chips - dry"
$ws.Range("C19").Value = 20.4
